$d = $word.ActiveDocument

$replacements = @(
    @("522÷8=65, 2", "600÷2=300, 0"),
    @("862÷5=172, 2", "380÷8=47, 4"),
    @("820÷3=273, 1", "463÷3=154, 1"),
    @("497÷9=55, 2", "446÷8=55, 6"),
    @("866÷7=123, 5", "295÷7=42, 1"),
    @("196÷6=32, 4", "824÷4=206, 0"),
    @("484÷4=121, 0", "406÷3=135, 1"),
    @("340÷9=37, 7", "451÷2=225, 1"),
    @("592÷2=296, 0", "293÷5=58, 3"),
    @("660÷6=110, 0", "198÷6=33, 0"),
    @("699÷5=139, 4", "431÷2=215, 1"),
    @("897÷8=112, 1", "919÷8=114, 7"),
    @("965÷9=107, 2", "237÷5=47, 2"),
    @("248÷6=41, 2", "558÷3=186, 0"),
    @("170÷6=28, 2", "851÷9=94, 5"),
    @("507÷4=126, 3", "376÷6=62, 4"),
    @("427÷4=106, 3", "968÷6=161, 2"),
    @("294÷6=49, 0", "164÷4=41, 0"),
    @("281÷2=140, 1", "849÷3=283, 0"),
    @("286÷4=71, 2", "900÷8=112, 4"),
    @("365÷2=182, 1", "509÷7=72, 5"),
    @("661÷3=220, 1", "526÷4=131, 2"),
    @("120÷9=13, 3", "727÷8=90, 7"),
    @("245÷5=49, 0", "981÷6=163, 3"),
    @("964÷2=482, 0", "418÷7=59, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
